# Re-run of the NATMI ligand-receptor pair analysis (Pdgfc-Pdgfra) after adding the
# "ECs" sending/target cluster to the dataset, per Dr Hou's advice. This replaces the
# previous 12-row (3x4 cluster) result table with the new 16-row (4x4 cluster) table
# and extends the sheet dimensions from A1:T13 to A1:T17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20
# row 2: ECs -> ECs
$data[0,0] = "ECs"
$data[0,1] = "Pdgfc"
$data[0,2] = "Pdgfra"
$data[0,3] = "ECs"
$data[0,4] = [double]"2"
$data[0,5] = [double]"0.6666666666666666"
$data[0,6] = [double]"0.1868766666666667"
$data[0,7] = [double]"0.56063"
$data[0,8] = [double]"0.01256665704529768"
$data[0,9] = [double]"0.01256665704529768"
$data[0,10] = [double]"3"
$data[0,11] = [double]"1"
$data[0,12] = [double]"3.535386"
$data[0,13] = [double]"10.606158"
$data[0,14] = [double]"0.01988747852527457"
$data[0,15] = [double]"0.01988747852527457"
$data[0,16] = [double]"0.66068115106"
$data[0,17] = [double]"5.94613035954"
$data[0,18] = [double]"0.0002499191221228481"
$data[0,19] = [double]"0.000249919122122848"
# row 3: ECs -> FAPs
$data[1,0] = "ECs"
$data[1,1] = "Pdgfc"
$data[1,2] = "Pdgfra"
$data[1,3] = "FAPs"
$data[1,4] = [double]"2"
$data[1,5] = [double]"0.6666666666666666"
$data[1,6] = [double]"0.1868766666666667"
$data[1,7] = [double]"0.56063"
$data[1,8] = [double]"0.01256665704529768"
$data[1,9] = [double]"0.01256665704529768"
$data[1,10] = [double]"3"
$data[1,11] = [double]"1"
$data[1,12] = [double]"173.8189136666666"
$data[1,13] = [double]"521.456741"
$data[1,14] = [double]"0.9777772251268709"
$data[1,15] = [double]"0.9777772251268707"
$data[1,16] = [double]"32.48269918964778"
$data[1,17] = [double]"292.3442927068299"
$data[1,18] = [double]"0.01228739105487221"
$data[1,19] = [double]"0.01228739105487221"
# row 4: ECs -> M2
$data[2,0] = "ECs"
$data[2,1] = "Pdgfc"
$data[2,2] = "Pdgfra"
$data[2,3] = "M2"
$data[2,4] = [double]"2"
$data[2,5] = [double]"0.6666666666666666"
$data[2,6] = [double]"0.1868766666666667"
$data[2,7] = [double]"0.56063"
$data[2,8] = [double]"0.01256665704529768"
$data[2,9] = [double]"0.01256665704529768"
$data[2,10] = [double]"3"
$data[2,11] = [double]"1"
$data[2,12] = [double]"0.06908833333333333"
$data[2,13] = [double]"0.207265"
$data[2,14] = [double]"0.000388640093475982"
$data[2,15] = [double]"0.0003886400934759819"
$data[2,16] = [double]"0.01291099743888889"
$data[2,17] = [double]"0.11619897695"
$data[2,18] = [double]"4.883906768765099e-06"
$data[2,19] = [double]"4.883906768765098e-06"
# row 5: ECs -> sCs
$data[3,0] = "ECs"
$data[3,1] = "Pdgfc"
$data[3,2] = "Pdgfra"
$data[3,3] = "sCs"
$data[3,4] = [double]"2"
$data[3,5] = [double]"0.6666666666666666"
$data[3,6] = [double]"0.1868766666666667"
$data[3,7] = [double]"0.56063"
$data[3,8] = [double]"0.01256665704529768"
$data[3,9] = [double]"0.01256665704529768"
$data[3,10] = [double]"3"
$data[3,11] = [double]"1"
$data[3,12] = [double]"0.346056"
$data[3,13] = [double]"1.038168"
$data[3,14] = [double]"0.001946656254378565"
$data[3,15] = [double]"0.001946656254378564"
$data[3,16] = [double]"0.06466979176"
$data[3,17] = [double]"0.58202812584"
$data[3,18] = [double]"2.446296153385919e-05"
$data[3,19] = [double]"2.446296153385918e-05"
# row 6: FAPs -> ECs
$data[4,0] = "FAPs"
$data[4,1] = "Pdgfc"
$data[4,2] = "Pdgfra"
$data[4,3] = "ECs"
$data[4,4] = [double]"3"
$data[4,5] = [double]"1"
$data[4,6] = [double]"3.160274"
$data[4,7] = [double]"9.480822"
$data[4,8] = [double]"0.212514918184031"
$data[4,9] = [double]"0.212514918184031"
$data[4,10] = [double]"3"
$data[4,11] = [double]"1"
$data[4,12] = [double]"3.535386"
$data[4,13] = [double]"10.606158"
$data[4,14] = [double]"0.01988747852527457"
$data[4,15] = [double]"0.01988747852527457"
$data[4,16] = [double]"11.172788455764"
$data[4,17] = [double]"100.555096101876"
$data[4,18] = [double]"0.004226385871685399"
$data[4,19] = [double]"0.004226385871685398"
# row 7: FAPs -> FAPs
$data[5,0] = "FAPs"
$data[5,1] = "Pdgfc"
$data[5,2] = "Pdgfra"
$data[5,3] = "FAPs"
$data[5,4] = [double]"3"
$data[5,5] = [double]"1"
$data[5,6] = [double]"3.160274"
$data[5,7] = [double]"9.480822"
$data[5,8] = [double]"0.212514918184031"
$data[5,9] = [double]"0.212514918184031"
$data[5,10] = [double]"3"
$data[5,11] = [double]"1"
$data[5,12] = [double]"173.8189136666666"
$data[5,13] = [double]"521.456741"
$data[5,14] = [double]"0.9777772251268709"
$data[5,15] = [double]"0.9777772251268707"
$data[5,16] = [double]"549.3153935690112"
$data[5,17] = [double]"4943.838542121101"
$data[5,18] = [double]"0.2077922470000458"
$data[5,19] = [double]"0.2077922470000458"
# row 8: FAPs -> M2
$data[6,0] = "FAPs"
$data[6,1] = "Pdgfc"
$data[6,2] = "Pdgfra"
$data[6,3] = "M2"
$data[6,4] = [double]"3"
$data[6,5] = [double]"1"
$data[6,6] = [double]"3.160274"
$data[6,7] = [double]"9.480822"
$data[6,8] = [double]"0.212514918184031"
$data[6,9] = [double]"0.212514918184031"
$data[6,10] = [double]"3"
$data[6,11] = [double]"1"
$data[6,12] = [double]"0.06908833333333333"
$data[6,13] = [double]"0.207265"
$data[6,14] = [double]"0.000388640093475982"
$data[6,15] = [double]"0.0003886400934759819"
$data[6,16] = [double]"0.2183380635366667"
$data[6,17] = [double]"1.96504257183"
$data[6,18] = [double]"8.259181766808246e-05"
$data[6,19] = [double]"8.259181766808244e-05"
# row 9: FAPs -> sCs
$data[7,0] = "FAPs"
$data[7,1] = "Pdgfc"
$data[7,2] = "Pdgfra"
$data[7,3] = "sCs"
$data[7,4] = [double]"3"
$data[7,5] = [double]"1"
$data[7,6] = [double]"3.160274"
$data[7,7] = [double]"9.480822"
$data[7,8] = [double]"0.212514918184031"
$data[7,9] = [double]"0.212514918184031"
$data[7,10] = [double]"3"
$data[7,11] = [double]"1"
$data[7,12] = [double]"0.346056"
$data[7,13] = [double]"1.038168"
$data[7,14] = [double]"0.001946656254378565"
$data[7,15] = [double]"0.001946656254378564"
$data[7,16] = [double]"1.093631779344"
$data[7,17] = [double]"9.842686014096"
$data[7,18] = [double]"0.0004136934946316929"
$data[7,19] = [double]"0.0004136934946316928"
# row 10: M2 -> ECs
$data[8,0] = "M2"
$data[8,1] = "Pdgfc"
$data[8,2] = "Pdgfra"
$data[8,3] = "ECs"
$data[8,4] = [double]"3"
$data[8,5] = [double]"1"
$data[8,6] = [double]"8.401240333333334"
$data[8,7] = [double]"25.203721"
$data[8,8] = [double]"0.5649475020465676"
$data[8,9] = [double]"0.5649475020465676"
$data[8,10] = [double]"3"
$data[8,11] = [double]"1"
$data[8,12] = [double]"3.535386"
$data[8,13] = [double]"10.606158"
$data[8,14] = [double]"0.01988747852527457"
$data[8,15] = [double]"0.01988747852527457"
$data[8,16] = [double]"29.70162745710201"
$data[8,17] = [double]"267.3146471139181"
$data[8,18] = [double]"0.01123538131485863"
$data[8,19] = [double]"0.01123538131485862"
# row 11: M2 -> FAPs
$data[9,0] = "M2"
$data[9,1] = "Pdgfc"
$data[9,2] = "Pdgfra"
$data[9,3] = "FAPs"
$data[9,4] = [double]"3"
$data[9,5] = [double]"1"
$data[9,6] = [double]"8.401240333333334"
$data[9,7] = [double]"25.203721"
$data[9,8] = [double]"0.5649475020465676"
$data[9,9] = [double]"0.5649475020465676"
$data[9,10] = [double]"3"
$data[9,11] = [double]"1"
$data[9,12] = [double]"173.8189136666666"
$data[9,13] = [double]"521.456741"
$data[9,14] = [double]"0.9777772251268709"
$data[9,15] = [double]"0.9777772251268707"
$data[9,16] = [double]"1460.294468192584"
$data[9,17] = [double]"13142.65021373326"
$data[9,18] = [double]"0.5523928008934501"
$data[9,19] = [double]"0.55239280089345"
# row 12: M2 -> M2
$data[10,0] = "M2"
$data[10,1] = "Pdgfc"
$data[10,2] = "Pdgfra"
$data[10,3] = "M2"
$data[10,4] = [double]"3"
$data[10,5] = [double]"1"
$data[10,6] = [double]"8.401240333333334"
$data[10,7] = [double]"25.203721"
$data[10,8] = [double]"0.5649475020465676"
$data[10,9] = [double]"0.5649475020465676"
$data[10,10] = [double]"3"
$data[10,11] = [double]"1"
$data[10,12] = [double]"0.06908833333333333"
$data[10,13] = [double]"0.207265"
$data[10,14] = [double]"0.000388640093475982"
$data[10,15] = [double]"0.0003886400934759819"
$data[10,16] = [double]"0.5804276925627778"
$data[10,17] = [double]"5.223849233065001"
$data[10,18] = [double]"0.0002195612500044006"
$data[10,19] = [double]"0.0002195612500044005"
# row 13: M2 -> sCs
$data[11,0] = "M2"
$data[11,1] = "Pdgfc"
$data[11,2] = "Pdgfra"
$data[11,3] = "sCs"
$data[11,4] = [double]"3"
$data[11,5] = [double]"1"
$data[11,6] = [double]"8.401240333333334"
$data[11,7] = [double]"25.203721"
$data[11,8] = [double]"0.5649475020465676"
$data[11,9] = [double]"0.5649475020465676"
$data[11,10] = [double]"3"
$data[11,11] = [double]"1"
$data[11,12] = [double]"0.346056"
$data[11,13] = [double]"1.038168"
$data[11,14] = [double]"0.001946656254378565"
$data[11,15] = [double]"0.001946656254378564"
$data[11,16] = [double]"2.907299624792"
$data[11,17] = [double]"26.165696623128"
$data[11,18] = [double]"0.001099758588254498"
$data[11,19] = [double]"0.001099758588254498"
# row 14: sCs -> ECs
$data[12,0] = "sCs"
$data[12,1] = "Pdgfc"
$data[12,2] = "Pdgfra"
$data[12,3] = "ECs"
$data[12,4] = [double]"3"
$data[12,5] = [double]"1"
$data[12,6] = [double]"3.122442666666667"
$data[12,7] = [double]"9.367328"
$data[12,8] = [double]"0.2099709227241037"
$data[12,9] = [double]"0.2099709227241037"
$data[12,10] = [double]"3"
$data[12,11] = [double]"1"
$data[12,12] = [double]"3.535386"
$data[12,13] = [double]"10.606158"
$data[12,14] = [double]"0.01988747852527457"
$data[12,15] = [double]"0.01988747852527457"
$data[12,16] = [double]"11.039040089536"
$data[12,17] = [double]"99.35136080582402"
$data[12,18] = [double]"0.0041757922166077"
$data[12,19] = [double]"0.004175792216607698"
# row 15: sCs -> FAPs
$data[13,0] = "sCs"
$data[13,1] = "Pdgfc"
$data[13,2] = "Pdgfra"
$data[13,3] = "FAPs"
$data[13,4] = [double]"3"
$data[13,5] = [double]"1"
$data[13,6] = [double]"3.122442666666667"
$data[13,7] = [double]"9.367328"
$data[13,8] = [double]"0.2099709227241037"
$data[13,9] = [double]"0.2099709227241037"
$data[13,10] = [double]"3"
$data[13,11] = [double]"1"
$data[13,12] = [double]"173.8189136666666"
$data[13,13] = [double]"521.456741"
$data[13,14] = [double]"0.9777772251268709"
$data[13,15] = [double]"0.9777772251268707"
$data[13,16] = [double]"542.7395923064497"
$data[13,17] = [double]"4884.656330758048"
$data[13,18] = [double]"0.2053047861785028"
$data[13,19] = [double]"0.2053047861785027"
# row 16: sCs -> M2
$data[14,0] = "sCs"
$data[14,1] = "Pdgfc"
$data[14,2] = "Pdgfra"
$data[14,3] = "M2"
$data[14,4] = [double]"3"
$data[14,5] = [double]"1"
$data[14,6] = [double]"3.122442666666667"
$data[14,7] = [double]"9.367328"
$data[14,8] = [double]"0.2099709227241037"
$data[14,9] = [double]"0.2099709227241037"
$data[14,10] = [double]"3"
$data[14,11] = [double]"1"
$data[14,12] = [double]"0.06908833333333333"
$data[14,13] = [double]"0.207265"
$data[14,14] = [double]"0.000388640093475982"
$data[14,15] = [double]"0.0003886400934759819"
$data[14,16] = [double]"0.2157243597688889"
$data[14,17] = [double]"1.94151923792"
$data[14,18] = [double]"8.160311903473386e-05"
$data[14,19] = [double]"8.160311903473384e-05"
# row 17: sCs -> sCs
$data[15,0] = "sCs"
$data[15,1] = "Pdgfc"
$data[15,2] = "Pdgfra"
$data[15,3] = "sCs"
$data[15,4] = [double]"3"
$data[15,5] = [double]"1"
$data[15,6] = [double]"3.122442666666667"
$data[15,7] = [double]"9.367328"
$data[15,8] = [double]"0.2099709227241037"
$data[15,9] = [double]"0.2099709227241037"
$data[15,10] = [double]"3"
$data[15,11] = [double]"1"
$data[15,12] = [double]"0.346056"
$data[15,13] = [double]"1.038168"
$data[15,14] = [double]"0.001946656254378565"
$data[15,15] = [double]"0.001946656254378564"
$data[15,16] = [double]"1.080540019456"
$data[15,17] = [double]"9.724860175104"
$data[15,18] = [double]"0.0004087412099585148"
$data[15,19] = [double]"0.0004087412099585147"

# Write the full A2:T17 block in one shot (row 1 header / styles are untouched)
$ws.Range("A2:T17").Value = $data
